$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 11429.75
$ws.Range("I9").Value = 13054.857
$ws.Range("J9").Value = 54
$ws.Range("K9").Value = 13054.857
$ws.Range("L9").Value = 54
$ws.Range("M9").Value = -12885.857
$ws.Range("N9").Value = -392
# Row 15
$ws.Range("H15").Value = 1427.0172
$ws.Range("I15").Value = 1427.0172
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4281.0516
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4112.0516
# Row 17
$ws.Range("H17").Value = 1222105
$ws.Range("I17").Value = 2357.875
$ws.Range("J17").Value = 2441852
$ws.Range("K17").Value = 7073.625
$ws.Range("L17").Value = 7325556
$ws.Range("M17").Value = -6905.625
$ws.Range("N17").Value = -7325892
# Row 76
$ws.Range("H76").Value = 5250
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5250
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 5250
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -5880
# Row 79
$ws.Range("H79").Value = 5250
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5250
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 5250
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -7434
# Row 137
$ws.Range("H137").Value = 2303.4119
$ws.Range("I137").Value = 2582.2173
$ws.Range("J137").Value = 1720.4546
$ws.Range("K137").Value = 7746.651899999999
$ws.Range("L137").Value = 5161.3638
$ws.Range("M137").Value = -5196.651899999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 500
$ws.Range("N11").Value = -788
# Row 32
$ws.Range("H32").Value = 7839.7534
$ws.Range("I32").Value = 3997.8833
$ws.Range("J32").Value = 33452.223
$ws.Range("K32").Value = 3997.8833
$ws.Range("L32").Value = 33452.223
$ws.Range("M32").Value = -3710.8833
$ws.Range("N32").Value = -34026.223
# Row 45
$ws.Range("H45").Value = 8011.25
$ws.Range("I45").Value = 10798.363
$ws.Range("J45").Value = 1879.6
$ws.Range("K45").Value = 10798.363
$ws.Range("L45").Value = 1879.6
$ws.Range("M45").Value = -10421.363
# Row 97
$ws.Range("H97").Value = 1067.6957
$ws.Range("I97").Value = 818.125
$ws.Range("J97").Value = 2731.5
$ws.Range("K97").Value = 818.125
$ws.Range("L97").Value = 2731.5
$ws.Range("M97").Value = -322.125

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5538.737
$ws.Range("I86").Value = 5983.5
$ws.Range("J86").Value = 3166.6667
$ws.Range("K86").Value = 5983.5
$ws.Range("L86").Value = 3166.6667
$ws.Range("M86").Value = -4860.5
$ws.Range("N86").Value = -5412.6667
# Row 89
$ws.Range("H89").Value = 5538.737
$ws.Range("I89").Value = 5983.5
$ws.Range("J89").Value = 3166.6667
$ws.Range("K89").Value = 29917.5
$ws.Range("L89").Value = 15833.3335
$ws.Range("M89").Value = -24301.5
$ws.Range("N89").Value = -27065.3335
# Row 94
$ws.Range("H94").Value = 897.1591
$ws.Range("I94").Value = 662.4706
$ws.Range("J94").Value = 1695.1
$ws.Range("K94").Value = 662.4706
$ws.Range("L94").Value = 1695.1
$ws.Range("M94").Value = -211.4706
$ws.Range("N94").Value = -2597.1
# Row 105
$ws.Range("H105").Value = 2477.8667
$ws.Range("I105").Value = 2663.625
$ws.Range("J105").Value = 2020.6154
$ws.Range("K105").Value = 2663.625
$ws.Range("L105").Value = 2020.6154
$ws.Range("M105").Value = -916.625
# Row 134
$ws.Range("H134").Value = 1961.7106
$ws.Range("I134").Value = 1492.6666
$ws.Range("J134").Value = 5057.4
$ws.Range("K134").Value = 4477.9998
$ws.Range("L134").Value = 15172.2
$ws.Range("M134").Value = -1942.9998

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 1287.5
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1287.5
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1287.5
$ws.Range("N14").Value = -1627.5
# Row 15
$ws.Range("H15").Value = 5997.5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 5997.5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5997.5
$ws.Range("N15").Value = -6337.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 124
$ws.Range("H124").Value = 30
$ws.Range("I124").Value = 30
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 90
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 4820
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
# Row 141
$ws.Range("H141").Value = 114842.555
$ws.Range("I141").Value = 1591.6666
$ws.Range("J141").Value = 341344.34
$ws.Range("K141").Value = 4774.9998
$ws.Range("L141").Value = 1024033.02
$ws.Range("M141").Value = 405.0002000000004
$ws.Range("N141").Value = -1034393.02

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 7942.231
$ws.Range("I2").Value = 11324.111
$ws.Range("J2").Value = 333
$ws.Range("K2").Value = 11324.111
$ws.Range("L2").Value = 333
$ws.Range("M2").Value = -11211.111
$ws.Range("N2").Value = -559
# Row 70
$ws.Range("H70").Value = 10959.632
$ws.Range("I70").Value = 5218.8667
$ws.Range("J70").Value = 32487.5
$ws.Range("K70").Value = 5218.8667
$ws.Range("L70").Value = 32487.5
$ws.Range("M70").Value = -4948.8667
# Row 73
$ws.Range("H73").Value = 10959.632
$ws.Range("I73").Value = 5218.8667
$ws.Range("J73").Value = 32487.5
$ws.Range("K73").Value = 5218.8667
$ws.Range("L73").Value = 32487.5
$ws.Range("M73").Value = -4282.8667
# Row 97
$ws.Range("H97").Value = 1291.5807
$ws.Range("I97").Value = 1412.2273
$ws.Range("J97").Value = 996.6667
$ws.Range("K97").Value = 1412.2273
$ws.Range("L97").Value = 996.6667
$ws.Range("M97").Value = -916.2273
# Row 113
$ws.Range("H113").Value = 3643.9524
$ws.Range("I113").Value = 3562
$ws.Range("J113").Value = 3992.25
$ws.Range("K113").Value = 3562
$ws.Range("L113").Value = 3992.25
$ws.Range("M113").Value = -1392
# Row 132
$ws.Range("H132").Value = 6026.32
$ws.Range("I132").Value = 2649.7144
$ws.Range("J132").Value = 23753.5
$ws.Range("K132").Value = 7949.1432
$ws.Range("L132").Value = 71260.5
$ws.Range("M132").Value = -5419.1432
$ws.Range("N132").Value = -76320.5
# Row 136
$ws.Range("H136").Value = 33059.08
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 33059.08
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 99177.24000000001
$ws.Range("N136").Value = -104277.24

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4540.7334
$ws.Range("I7").Value = 2783.8823
$ws.Range("J7").Value = 6838.154
$ws.Range("K7").Value = 2783.8823
$ws.Range("L7").Value = 6838.154
$ws.Range("M7").Value = -2671.8823
$ws.Range("N7").Value = -7062.154
# Row 16
$ws.Range("H16").Value = 2724.3928
$ws.Range("I16").Value = 2263.818
$ws.Range("J16").Value = 4413.1665
$ws.Range("K16").Value = 2263.818
$ws.Range("L16").Value = 4413.1665
$ws.Range("M16").Value = -2093.818
# Row 22
$ws.Range("H22").Value = 2455.4092
$ws.Range("I22").Value = 1388.4286
$ws.Range("J22").Value = 2953.3333
$ws.Range("K22").Value = 1388.4286
$ws.Range("L22").Value = 2953.3333
$ws.Range("M22").Value = -1093.4286
$ws.Range("N22").Value = -3543.3333
# Row 27
$ws.Range("H27").Value = 2455.4092
$ws.Range("I27").Value = 1388.4286
$ws.Range("J27").Value = 2953.3333
$ws.Range("K27").Value = 1388.4286
$ws.Range("L27").Value = 2953.3333
$ws.Range("M27").Value = -1281.4286
$ws.Range("N27").Value = -3167.3333
# Row 61
$ws.Range("H61").Value = 2089.6667
$ws.Range("I61").Value = 1955.84
$ws.Range("J61").Value = 3762.5
$ws.Range("K61").Value = 1955.84
$ws.Range("L61").Value = 3762.5
$ws.Range("M61").Value = -1753.84
# Row 87
$ws.Range("H87").Value = 104000
$ws.Range("I87").Value = 9000
$ws.Range("J87").Value = 199000
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 199000
$ws.Range("M87").Value = -7877
$ws.Range("N87").Value = -201246
# Row 90
$ws.Range("H90").Value = 104000
$ws.Range("I90").Value = 9000
$ws.Range("J90").Value = 199000
$ws.Range("K90").Value = 27000
$ws.Range("L90").Value = 597000
$ws.Range("M90").Value = -21384
$ws.Range("N90").Value = -608232
# Row 113
$ws.Range("H113").Value = 2089.6667
$ws.Range("I113").Value = 1955.84
$ws.Range("J113").Value = 3762.5
$ws.Range("K113").Value = 1955.84
$ws.Range("L113").Value = 3762.5
$ws.Range("M113").Value = 214.1600000000001
# Row 126
$ws.Range("H126").Value = 4540.7334
$ws.Range("I126").Value = 2783.8823
$ws.Range("J126").Value = 6838.154
$ws.Range("K126").Value = 8351.6469
$ws.Range("L126").Value = 20514.462
$ws.Range("M126").Value = -5881.6469
$ws.Range("N126").Value = -25454.462
# Row 136
$ws.Range("H136").Value = 3728.0952
$ws.Range("I136").Value = 3146.842
$ws.Range("J136").Value = 9250
$ws.Range("K136").Value = 9440.526
$ws.Range("L136").Value = 27750
$ws.Range("M136").Value = -6890.526

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 28810.166
$ws.Range("I107").Value = 1076.5667
$ws.Range("J107").Value = 167478.17
$ws.Range("K107").Value = 3229.7001
$ws.Range("L107").Value = 502434.51
$ws.Range("M107").Value = -1309.7001
# Row 113
$ws.Range("H113").Value = 553.9459000000001
$ws.Range("I113").Value = 540.25
$ws.Range("J113").Value = 596.55554
$ws.Range("K113").Value = 1620.75
$ws.Range("L113").Value = 1789.66662
$ws.Range("M113").Value = 549.25
$ws.Range("N113").Value = -6129.66662
# Row 138
$ws.Range("H138").Value = 105000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 105000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 105000
$ws.Range("N138").Value = -115280

